$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Retangulo 15" - split the single-paragraph "Home3" run into two
# paragraphs: "Home" and "3" (a carriage return inserts a new a:p that
# inherits the existing paragraph/run formatting).
$s.Shapes.Item(11).TextFrame.TextRange.Text = "Home`r3"

# "Retangulo 17" - simple text fix: "Home3" -> "Homer3".
$s.Shapes.Item(12).TextFrame.TextRange.Text = "Homer3"
